# Insert a new weekly data row at row 104 (shifts existing rows 104-170
# down to 105-171), and populate the new row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 104:170 down by one row, leaving row 104 empty for the new entry.
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with the new weekly record.
$ws.Cells.Item(104, 1).Value = 7
$ws.Cells.Item(104, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(104, 3).Value = "Ñuble"
$ws.Cells.Item(104, 4).Value = 44488
$ws.Cells.Item(104, 5).Value = 16
$ws.Cells.Item(104, 6).Value = 100112008
$ws.Cells.Item(104, 7).Value = "Coliflor"
$ws.Cells.Item(104, 8).Value = "Sin especificar"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 180
$ws.Cells.Item(104, 11).Value = 700
$ws.Cells.Item(104, 12).Value = 800
$ws.Cells.Item(104, 13).Value = 750
$ws.Cells.Item(104, 14).Value = '$/unidad'
$ws.Cells.Item(104, 15).Value = "Región del Maule"
$ws.Cells.Item(104, 16).Value = 750
$ws.Cells.Item(104, 17).Value = 1
$ws.Cells.Item(104, 18).Value = "Hortaliza"
